# Apply strikethrough formatting to the "Digital Portfolio" intro paragraphs
# and the deliverable description, and drop the stray _GoBack bookmark —
# matching the target revision of the document.

$d = $word.ActiveDocument

$targets = @(
    "A digital portfolio provides a platform to highlight your",
    "This portfolio will serve as a digital representation of your academic",
    "Digital Portfolio:",
    "Students must provide a fully accessible digital portfolio"
)

foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text
    foreach ($prefix in $targets) {
        if ($t.StartsWith($prefix)) {
            $p.Range.Font.StrikeThrough = $true
            break
        }
    }
}

# Remove the leftover "_GoBack" bookmark (present at the end of the
# "Platform Choice" bullet) which Word no longer keeps on a clean re-save.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}
